$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Expenses sheet: E10 and E11 change from numeric storage to text storage
#    (same displayed values: 700 and 48), matching the diff's
#    `<c r="E10">` -> `<c r="E10" t="str">` change.
# ---------------------------------------------------------------------------
$wsExpenses = $wb.Worksheets.Item("Expenses")

$wsExpenses.Cells.Item(10, 5).NumberFormat = "@"
$wsExpenses.Cells.Item(10, 5).Value = "700"
$wsExpenses.Cells.Item(10, 5).ClearFormats()

$wsExpenses.Cells.Item(11, 5).NumberFormat = "@"
$wsExpenses.Cells.Item(11, 5).Value = "48"
$wsExpenses.Cells.Item(11, 5).ClearFormats()

# ---------------------------------------------------------------------------
# 2) Users sheet: add permissions column (F) values and refresh modifiedAt (H)
#    timestamps for rows 2-4.
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")

$wsUsers.Cells.Item(2, 6).Value = "read,write,delete,admin"
$wsUsers.Cells.Item(2, 8).Value = "2025-04-29T04:58:50.090Z"

$wsUsers.Cells.Item(3, 6).Value = "read,write,delete"
$wsUsers.Cells.Item(3, 8).Value = "2025-04-29T04:59:16.302Z"

$wsUsers.Cells.Item(4, 6).Value = "read"
$wsUsers.Cells.Item(4, 8).Value = "2025-04-29T04:59:22.458Z"

# ---------------------------------------------------------------------------
# 3) AuditLog sheet: append 6 new audit rows (18-23) recording the
#    permission updates made to the three users above.
# ---------------------------------------------------------------------------
$wsAudit = $wb.Worksheets.Item("AuditLog")

# NOTE: this runtime's PowerShell parser doesn't handle named (-Flag value)
# arguments for user-defined functions reliably, so Set-AuditRow is called
# with plain positional arguments.
function Set-AuditRow($Row, $Id, $EntityType, $EntityId, $Action, $UserId, $Username, $Timestamp, $Changes, $Description) {
    $wsAudit.Cells.Item($Row, 1).Value = $Id
    $wsAudit.Cells.Item($Row, 2).Value = $EntityType

    $wsAudit.Cells.Item($Row, 3).NumberFormat = "@"
    $wsAudit.Cells.Item($Row, 3).Value = $EntityId
    $wsAudit.Cells.Item($Row, 3).ClearFormats()

    $wsAudit.Cells.Item($Row, 4).Value = $Action

    $wsAudit.Cells.Item($Row, 5).NumberFormat = "@"
    $wsAudit.Cells.Item($Row, 5).Value = $UserId
    $wsAudit.Cells.Item($Row, 5).ClearFormats()

    $wsAudit.Cells.Item($Row, 6).Value = $Username
    $wsAudit.Cells.Item($Row, 7).Value = $Timestamp
    $wsAudit.Cells.Item($Row, 8).Value = $Changes
    $wsAudit.Cells.Item($Row, 9).Value = $Description
}

$changes1 = '{"before":{"id":"1","username":"admin","name":"Lisa Williams","email":"admin@example.com","role":"admin","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-24T00:54:41.829Z","status":"active"},"after":{"id":"1","username":"admin","name":"Lisa Williams","email":"admin@example.com","role":"admin","permissions":["read","write","delete","admin"],"modifiedBy":"admin","modifiedAt":"2025-04-29T04:58:50.090Z","status":"active"}}'

$changes2 = '{"before":{"id":"2","username":"user","name":"Regular User","email":"user@example.com","role":"user","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-28T23:24:30.933Z","status":"active"},"after":{"id":"2","username":"user","name":"Regular User","email":"user@example.com","role":"user","permissions":["read","write","delete"],"modifiedBy":"admin","modifiedAt":"2025-04-29T04:59:16.302Z","status":"active"}}'

$changes3 = '{"before":{"id":"3","username":"viewer","name":"Viewer","email":"viewer@example.com","role":"viewer","permissions":[],"modifiedBy":"admin","modifiedAt":"2025-04-28T23:24:41.205Z","status":"active"},"after":{"id":"3","username":"viewer","name":"Viewer","email":"viewer@example.com","role":"viewer","permissions":["read"],"modifiedBy":"admin","modifiedAt":"2025-04-29T04:59:22.458Z","status":"active"}}'

Set-AuditRow 18 "AUDIT1745902730091" "Users" "1" "UPDATE" "1" "admin" "2025-04-29T04:58:50.091Z" $changes1 "Updated User 1"
Set-AuditRow 19 "AUDIT1745902730091" "Users" "1" "UPDATE" "1" "admin" "2025-04-29T04:58:50.091Z" $changes1 "Updated User 1"

Set-AuditRow 20 "AUDIT1745902756302" "Users" "2" "UPDATE" "1" "admin" "2025-04-29T04:59:16.302Z" $changes2 "Updated User 2"
Set-AuditRow 21 "AUDIT1745902756302" "Users" "2" "UPDATE" "1" "admin" "2025-04-29T04:59:16.302Z" $changes2 "Updated User 2"

Set-AuditRow 22 "AUDIT1745902762458" "Users" "3" "UPDATE" "1" "admin" "2025-04-29T04:59:22.458Z" $changes3 "Updated User 3"
Set-AuditRow 23 "AUDIT1745902762458" "Users" "3" "UPDATE" "1" "admin" "2025-04-29T04:59:22.458Z" $changes3 "Updated User 3"
